$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused rows (7:23), which also shrinks the used range/dimension to A1:A6
$ws.Range("A7:A23").ClearContents()

# Consolidate each card's fields into a single Python-tuple-like string per row
$ws.Range("A2").Value = '(''Arc Lightning'', [''{2}{R}'', ''Sorcery'', ''Arc Lightning deals 3 damage divided as you choose among one, two, or three targets.''])'
$ws.Range("A3").Value = '(''Dauthi Slayer'', [''{B}{B}'', ''Creature — Dauthi Soldier'', ''Shadow (This creature can block or be blocked by only creatures with shadow.)'', ''Dauthi Slayer attacks each combat if able.'', ''2/2''])'
$ws.Range("A4").Value = '(''Island'', [''Basic Land — Island'', ''({T}: Add {U}.)''])'
$ws.Range("A5").Value = '(''Mana Leak'', [''{1}{U}'', ''Instant'', ''Counter target spell unless its controller pays {3}.''])'
$ws.Range("A6").Value = '("Man-o''-War", [''{2}{U}'', ''Creature — Jellyfish'', ''When Man-o’-War enters the battlefield, return target creature to its owner’s hand.'', ''2/2''])'
